# Apply edits to Sheet1 of the Game.xlsx workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 3 values
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = -6
$ws.Range("K3").Value = 10

# Row 5 values
$ws.Range("G5").Value = 5
$ws.Range("I5").Value = 7
$ws.Range("K5").Value = 7

# Row 7 values
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 8
$ws.Range("K7").Value = 1

# Update the selected cell to match the saved view state
$ws.Range("N17").Select()
